$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header values (row 1, columns B-E)
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Update row 2 values (columns B-E)
$ws.Range("B2").Value = 39.737113767944933
$ws.Range("C2").Value = 39.71845511051913
$ws.Range("D2").Value = 42.885043917306568
$ws.Range("E2").Value = 35.739962095558155

# Update row 3 values (columns B-E)
$ws.Range("B3").Value = 57.342222431918422
$ws.Range("C3").Value = 43.174736418035501
$ws.Range("D3").Value = 45.329749470807954
$ws.Range("E3").Value = 17.073207678383952

# Update selection to match the new active range
$ws.Range("B1:E3").Select()
